$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: combined "timestamp{-}name" headers -> uppercase the name portion
$ws.Range("A1").Value = "22.09.2021 13:28:19{-}А КТО КТО КТО"
$ws.Range("B1").Value = "22.09.2021 13:29:13{-}ФАМИЛИЯ ИМЯ"
$ws.Range("C1").Value = "22.09.2021 13:47:58{-}ДЖЛДЛООГРИДЛ"
$ws.Range("D1").Value = "22.09.2021 13:48:39{-}ОРПРДЛО"
$ws.Range("E1").Value = "22.09.2021 13:49:01{-}ПЕТЯ"
$ws.Range("F1").Value = "22.09.2021 13:50:34{-}ВАСЯ"
$ws.Range("G1").Value = "28.09.2021 10:43:58{-}ОЫВАЖЫВА"
$ws.Range("H1").Value = "28.09.2021 14:57:51{-}SDFJDS"
$ws.Range("I1").Value = "28.09.2021 14:59:57{-}KDFJGDFHGFGHJF"
$ws.Range("J1").Value = "05.10.2021 09:12:37{-}Я"
$ws.Range("K1").Value = "05.10.2021 09:13:01{-}SDFKLSDFSDGSD"
$ws.Range("L1").Value = "06.10.2021 15:01:45{-}ФАМИЛИЯ ИМЯ"
$ws.Range("M1").Value = "12.10.2021 17:10:17{-}LKJG;LDFJGILJG"
$ws.Range("N1").Value = "12.10.2021 17:10:46{-}LKJG;LDFJGILJG"
$ws.Range("O1").Value = "13.10.2021 09:21:40{-}ОЗЕРОВА ЕЛИЗАВЕТА ДЕНИСОВНА"
$ws.Range("P1").Value = "13.10.2021 09:28:48{-}ЛОСКУТОВА"
$ws.Range("Q1").Value = "13.10.2021 09:34:31{-}ТУМАНОВА ИРИНА ЛЕОНИДОВНА"
$ws.Range("R1").Value = "13.10.2021 09:39:34{-}РОМАНЦОВ ДМИТРИЙ ИЛЬИЧ"
$ws.Range("S1").Value = "13.10.2021 09:39:40{-}АНДРЕЕВА СНЕЖАНА ГЕНАДЬЕВНА"
$ws.Range("T1").Value = "13.10.2021 09:44:54{-}БЕЛОУСОВА НАТАЛЬЯ"
$ws.Range("U1").Value = "13.10.2021 09:51:08{-}ПЕТРОВА МАРИНА АЛЕКСЕЕВНА"
$ws.Range("V1").Value = "20.10.2021 09:09:04{-}ОЗЕРОВА ЕЛИЗАВЕТА ДЕНИСОВНА"
$ws.Range("W1").Value = "20.10.2021 09:17:40{-}КЕКИНА ЛЮДМИЛА ФЕДОРОВНА"
$ws.Range("X1").Value = "20.10.2021 09:48:00{-}РОМАНЦОВ ДМИТРИЙ ИЛЬИЧ"
$ws.Range("Y1").Value = "10.01.2022 17:28:01{-}ФАМИЛИЯ ИМЯ"
$ws.Range("Z1").Value = "19.01.2022 10:22:18{-}ВАЖЕНИНА"
$ws.Range("AA1").Value = "19.01.2022 14:00:29{-}ФАМИЛИЯ ИМЯ ТЕСТИРУЕМ 19.01.2022"

# Row 7: "Имя пользователя (тестируемого)" plain name values -> uppercase
$ws.Range("A7").Value = "А КТО КТО КТО"
$ws.Range("B7").Value = "ФАМИЛИЯ ИМЯ"
$ws.Range("C7").Value = "ДЖЛДЛООГРИДЛ"
$ws.Range("D7").Value = "ОРПРДЛО"
$ws.Range("E7").Value = "ПЕТЯ"
$ws.Range("F7").Value = "ВАСЯ"
$ws.Range("G7").Value = "ОЫВАЖЫВА"
$ws.Range("H7").Value = "SDFJDS"
$ws.Range("I7").Value = "KDFJGDFHGFGHJF"
$ws.Range("J7").Value = "Я"
$ws.Range("K7").Value = "SDFKLSDFSDGSD"
$ws.Range("L7").Value = "ФАМИЛИЯ ИМЯ"
$ws.Range("M7").Value = "LKJG;LDFJGILJG"
$ws.Range("N7").Value = "LKJG;LDFJGILJG"
$ws.Range("O7").Value = "ОЗЕРОВА ЕЛИЗАВЕТА ДЕНИСОВНА"
$ws.Range("P7").Value = "ЛОСКУТОВА"
$ws.Range("Q7").Value = "ТУМАНОВА ИРИНА ЛЕОНИДОВНА"
$ws.Range("R7").Value = "РОМАНЦОВ ДМИТРИЙ ИЛЬИЧ"
$ws.Range("S7").Value = "АНДРЕЕВА СНЕЖАНА ГЕНАДЬЕВНА"
$ws.Range("T7").Value = "БЕЛОУСОВА НАТАЛЬЯ"
$ws.Range("U7").Value = "ПЕТРОВА МАРИНА АЛЕКСЕЕВНА"
$ws.Range("V7").Value = "ОЗЕРОВА ЕЛИЗАВЕТА ДЕНИСОВНА"
$ws.Range("W7").Value = "КЕКИНА ЛЮДМИЛА ФЕДОРОВНА"
$ws.Range("X7").Value = "РОМАНЦОВ ДМИТРИЙ ИЛЬИЧ"
$ws.Range("Y7").Value = "ФАМИЛИЯ ИМЯ"
$ws.Range("Z7").Value = "ВАЖЕНИНА"
$ws.Range("AA7").Value = "ФАМИЛИЯ ИМЯ ТЕСТИРУЕМ 19.01.2022"
